$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.181928873062134
$ws.Range("B1").Value = 2.329928159713745
$ws.Range("C1").Value = 3.831974983215332
$ws.Range("D1").Value = 3.096009492874146
$ws.Range("E1").Value = 1.142207145690918
